# Insert a new data row at row 147 (pushing the existing rows 147:219 down
# to 148:220) and populate the new row with the latest price record.
# This matches the diff: dimension grows from A1:R219 to A1:R220 and every
# existing record from row 147 onward shifts down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("147:147").Insert()

$ws.Range("A147").Value2 = 10
$ws.Range("B147").Value2 = "Vega Modelo de Temuco"
$ws.Range("C147").Value2 = "La Araucanía"
$ws.Range("D147").Value2 = 44460
$ws.Range("E147").Value2 = 9
$ws.Range("F147").Value2 = 100112040
$ws.Range("G147").Value2 = "Cilantro"
$ws.Range("H147").Value2 = "Sin especificar"
$ws.Range("I147").Value2 = "Primera"
$ws.Range("J147").Value2 = 30
$ws.Range("K147").Value2 = 6000
$ws.Range("L147").Value2 = 7000
$ws.Range("M147").Value2 = 6333
$ws.Range("N147").Value2 = "$/docena de atados (2 kilos)"
$ws.Range("O147").Value2 = "Provincia de Cautín"
$ws.Range("P147").Value2 = 3166
$ws.Range("Q147").Value2 = 2
$ws.Range("R147").Value2 = "Hortaliza"
